$d = $word.ActiveDocument

$searchText = "ушли под воду. "
$replacementText = "ушли под воду. Для первых, никогда не блиставших архитектурными шедеврами городов и не привязанных к одному месту, это не стало большим испытанием, и, хотя стычки между племенами участились, мудрость вождей и угруннов не дала разразиться гражданской войне.`rДля кадалов же это было очередное испытание веры и стойкости. С островами были утрачены как многие технологии, так и мастера, владевшие ими. Цивилизация ящеров откатилась на несколько тысяч лет назад: у них осталось лишь примитивные орудия и механизмы. Кадалы всегда находились под покровительством Айгне, но по иронии судьбы они же больше всего страдали от его действий и последствий этих действий. Но как не дано рыбе летать, так и простому обывателю не дано постичь, терзалось этим само божество.`rХотя Дитя Света родилось из осколков Айгне, оно обладало собственным сознание и собственной волей. Узрев Стигию и разочаровавшись, оно решило создать на ней новую цивилизацию, могущественней и крепче прежних, защищенную от невзгод и печалей. Но перед этим нужно было подготовить сам мир к творению.`rВ первую очередь Дитя Света создало защитников для новой расы – Светлое Воинство ангелов, воинов недюжинной силы и красоты, безукоризненно верных своему создателю. В стремлении не допустить предательства, Дитя Света соединило разумы ангелов в один поток, так что каждый знал мысли каждого, хотя и мог действовать самостоятельно.`rВидя жизнь Стигии, Дитя Света приняло решение создать ангелов по тем же принципам. Они получили плоть и кровь, потребность в еде и воде, способность к продолжению рода. Но помня об их истинном назначении божество лишило их любых эмоций, кроме любви друг к другу, к своему творцу и к будущей расе и гнева ко всему, что обратилось бы против созданий Дитя Света. `rПервым заданием Светлого Воинства стало очищение самых благоприятных территорий континента, лежавшего севернее Дунхайма, от оскверненных источников и их порождений. Сейчас там располагаются королевство людей Нимлис, земли нескольких Кланов дворфов и королевство эльфов Ха-Ли-Най. Увы, ангелы не были способны изгнать хаос из недр Стигии, потому ценой многих жертв, им удалось лишь запечатать порченные источники, перебив все, что вышло из них.`rПосле этого Дитя Света размножило растения и животных, которые стали бы едой и помощью для будущей расы. В лесах стало больше дичи, а в реках – рыбы, землю устлали ягоды, плодоносящие деревья стали давать больше фруктов. Убедившись, что хищники и добыча находятся в равновесии, земля плодоносит, Дитя Света сотворило Древних.`rОбликом своим они были похожи на ангелов, однако, были лишены крыльев и лишь едва-едва обладали способностью использовать силу Света. Но даже тех возможностей, что были им даны, Древним – или, как они были названы Дитем Света, атлантам, - хватало на творение невероятных чудес: от простых магических огоньков до сложных, осязаемых иллюзий и управления погодой.`rУ атлантов не было врагов, потому Дитя Света дало им лишь орудия для охоты, добычи и созидания. До сих пор в отдаленных уголках Северного Континента под землей можно найти остатки их городов – изящные, прекрасные, с невероятными фресками и декоративными элементами. Древние не знали горестей: живя под защитой ангелов, они имели возможность творить и развиваться. Поговаривают, что многие громадные механизмы гномов остались еще от атлантов.`rИх численность быстро росла, и Древние начали расселяться по Северному Континенту. Первыми ушли те, что потом станут эльфами. Они двинулись на восток, в теплые леса и поля, к озерам и речушкам, влекомые красотой этих мест. Найдя эти края едва ли не раем на земле, они продолжили совершенствовать искусства и ремесла, заселив весь восток континента.`rПосле них ушли те, что станут дворфами. Они пошли на север, желая изучить край гор, непроходимых лесов и холода. Найдя в недрах гор богатые жилы металлов и драгоценные камни, они решили поселиться под землей. Сейчас в мире нет, пожалуй, ни одной горы, где бы не побывал хотя бы один дворф.`rТретьими ушли те, что станут варкарами. Вопреки запретам Дитя Света и ангелов, они двинулись на запад, туда, где оскверненные источники оставались открытыми. Судьба их печальна: близость хаоса привлекла на эти земли взор Дитя Тьмы. Варкары пали под его могуществом, исказились, стали двуполыми уродливыми пародиями на своих могучих предков-атлантов. Впрочем, с их точки зрения, все сложилось более чем удачно.`rТе, кто остались мало-по-малу заселили все оставшиеся земли континента, от ледяных северных пустошей до жаркого побережья Великого Моря, разделяющего Дунхайм и Северный Континент. Они могли бы сохранить свое естество, подобно эльфам, но этому не суждено было сбыться.`rПечати, что были оставлены ангелами на оскрверненных источниках, начали рушиться под напором рвущейся наружу материи хаоса, и Дитя Тьмы, до того времени занимавшееся лишь варкарами, играясь с ними и уродуя их все больше, подобно безумному и жестокому ребенку, обратило свой взор на восток. `r"

$d.Content.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, $replacementText, 2)

$lastPara = $d.Paragraphs.Last
$prevPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$markRange = $d.Range($prevPara.Range.End - 1, $prevPara.Range.End)
$markRange.Delete()
